$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# 1) "W" + bookmarkStart/_GoBack/bookmarkEnd + "hat have you done?" collapse
#    into a single run reading "What have you done?" (the bookmark is
#    removed from here; it reappears on the final paragraph in step 3).
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstRange = $firstPara.Range
$firstSelection = $d.Range($firstRange.Start, $firstRange.End - 1)

$mergedRunXml = '<w:p xmlns:w="' + $wNs + '">' + `
    '<w:r w:rsidRPr="006127BB"><w:rPr><w:b/><w:i w:val="0"/><w:sz w:val="21"/><w:szCs w:val="22"/></w:rPr><w:t>What have you done?</w:t></w:r>' + `
    '</w:p>'

$firstSelection.InsertXML($mergedRunXml)

# ---------------------------------------------------------------------------
# 2) Append a new sentence as its own run at the end of the closing answer
#    paragraph ("Through use of time management..."), then add a bookmark
#    (_GoBack) inside the final, now-empty trailing paragraph.
# ---------------------------------------------------------------------------

# Locate the paragraph that ends with "...working together. " - it is the
# second-to-last paragraph in the body (the very last paragraph is empty).
$targetIndex = $d.Paragraphs.Count - 1
$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $targetPara.Range

# Insert a new paragraph (carrying the ORIGINAL paragraph's own identity/
# properties) right after the target paragraph's text, containing the new
# run. Inserting xml wrapped in <w:p> at the collapsed point immediately
# before the paragraph mark splices in a fresh sibling paragraph.
$insertPos = $d.Range($targetRange.End - 1, $targetRange.End - 1)

$newParaXml = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="24E64DF8" w14:textId="77777777" w:rsidR="00A840C2" w:rsidRPr="006127BB" w:rsidRDefault="00A840C2" w:rsidP="00A840C2">' + `
    '<w:pPr><w:spacing w:before="120" w:after="360" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:right="360"/><w:rPr><w:i w:val="0"/><w:sz w:val="21"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:i w:val="0"/><w:sz w:val="21"/><w:szCs w:val="22"/></w:rPr><w:t>The thoughts, and opinions of others can shed light upon something I may have been previously unaware of.</w:t></w:r>' + `
    '</w:p>'

$insertPos.InsertXML($newParaXml)

# Merge the original paragraph with the freshly inserted one by deleting the
# paragraph mark between them - this keeps both runs distinct while folding
# them into a single <w:p> (the second paragraph's identity/pPr - which we
# set above to match the original - becomes the identity of the merged one).
$mergeEnd = $d.Paragraphs.Item($targetIndex).Range.End
$mark = $d.Range($mergeEnd - 1, $mergeEnd)
$mark.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) Add the _GoBack bookmark into the final (empty) trailing paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastStart = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$bookmarkXml = '<w:p xmlns:w="' + $wNs + '"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$lastStart.InsertXML($bookmarkXml)
